# Add a new event (Id 18, "Студенческая олимпиада" / "Лингвистика", 20.12.2020)
# with three volunteer registrations — two existing volunteers (3 and 5) and
# one brand-new volunteer (4, Кудлай Полина Александровна). This mirrors the
# "someevents / userregisteronevent / alluser" reporting rows being appended
# to the registration log, plus a trailing orphan date value in A10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: event 18, volunteer 3 (Суханов Игорь Константинович) ---
$ws.Range("A7").Value = 18
$ws.Range("B7").Value = "Студенческая олимпиада"
$ws.Range("C7").Value = "Лингвистика"
$ws.Range("D7").Value = "20.12.2020"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = "Суханов"
$ws.Range("G7").Value = "Игорь"
$ws.Range("H7").Value = "Константинович"
$ws.Range("J7").Value = "you@me.they"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "456542145"
$ws.Range("L7").Value = "21.11.2001"
$ws.Range("M7").Value = 1

# --- Row 8: event 18, volunteer 4 (new volunteer: Кудлай Полина Александровна) ---
$ws.Range("A8").Value = 18
$ws.Range("B8").Value = "Студенческая олимпиада"
$ws.Range("C8").Value = "Лингвистика"
$ws.Range("D8").Value = "20.12.2020"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Кудлай"
$ws.Range("G8").Value = "Полина"
$ws.Range("H8").Value = "Александровна"
$ws.Range("I8").Value = "Иноватика"
$ws.Range("J8").Value = "pol@ina.com"
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "789654"
$ws.Range("L8").Value = "15.11.1999"
$ws.Range("M8").Value = 0

# --- Row 9: event 18, volunteer 5 (Мулькин Певел Сергеевич) ---
$ws.Range("A9").Value = 18
$ws.Range("B9").Value = "Студенческая олимпиада"
$ws.Range("C9").Value = "Лингвистика"
$ws.Range("D9").Value = "20.12.2020"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = "Мулькин"
$ws.Range("G9").Value = "Певел"
$ws.Range("H9").Value = "Сергеевич"
$ws.Range("I9").Value = "ЮрФУ"
$ws.Range("J9").Value = "patro1@yandex.ru"
$ws.Range("K9").Value = "223-322-223"
# "11.11.2011" is day<=12, which Excel's auto-detect reads as an ambiguous
# date rather than literal text, unlike the other DD.MM.YYYY strings here -
# force it to remain text to match the other date-string cells in this sheet.
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "11.11.2011"
$ws.Range("M9").Value = 0

# --- Row 10: stray reformatted date value (was A7 before the new rows pushed it down) ---
$ws.Range("A10").Value = "2019-07-24 00:00:00"
